$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Columns.Item(13).Delete()
$excel.ActiveWindow.Zoom = 100
